$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table with newly scraped values.
# Column D ("Price") holds numeric-looking text (e.g. "4.00", "520.07",
# "4.034.28", "0.0₃0862") that must stay literal text. Forcing the cell to
# Text format ("@") before assigning the value stops Excel from silently
# re-interpreting it as a floating point number (which would corrupt
# trailing zeros / introduce binary rounding noise). The format is then
# restored to the default "Normal" style so no stray cell styling remains.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.192.46"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.034.28"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.07"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.98"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.727"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +18.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.026.50"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.51%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.779"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.31%  "
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.94"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.07"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.678.71"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.045.74"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.24"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.17"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.22"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.105.19"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "444.14"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.18"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.57"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.90"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.49"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.02"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.74"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.82"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.23"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.68"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("E33").Value = "  +2.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "676.07"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.75"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "66.75"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.24"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0862"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.426"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.52"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.152"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0500"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("E46").Value = "  +14.23%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.50"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.21"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.81%  "
$ws.Range("E51").Value = "  +2.81%  "
